# Refresh the cryptocurrency price/volume snapshot (GitHub Actions scrape update).
# Source data is stored as plain text (coinranking.com scrape), so numeric-looking
# values are written with a leading apostrophe to keep Excel from coercing them to
# the Number type (matches the original inlineStr/text cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '70.102.43'
$ws.Cells.Item(2, 5).Value = '  -1.24%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.782.10'
$ws.Cells.Item(3, 5).Value = '  +2.91%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.00%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''618.14'
$ws.Cells.Item(5, 5).Value = '  +3.26%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''178.26'
$ws.Cells.Item(6, 5).Value = '  -3.55%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '3.779.50'
$ws.Cells.Item(7, 5).Value = '  +2.92%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.01%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.535'
$ws.Cells.Item(9, 5).Value = '  -0.19%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''0.172'
$ws.Cells.Item(10, 5).Value = '  +5.86%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''6.36'
$ws.Cells.Item(11, 5).Value = '  -2.84%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''0.492'
$ws.Cells.Item(12, 5).Value = '  -1.67%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''41.10'
$ws.Cells.Item(13, 5).Value = '  +2.90%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''0.0000264'
$ws.Cells.Item(14, 5).Value = '  +4.11%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '4.415.04'
$ws.Cells.Item(15, 5).Value = '  +3.01%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '3.782.33'
$ws.Cells.Item(16, 5).Value = '  +3.08%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '70.145.26'
$ws.Cells.Item(17, 5).Value = '  -1.20%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  -0.05%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +1.57%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''16.85'
$ws.Cells.Item(20, 5).Value = '  -1.59%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''513.04'
$ws.Cells.Item(21, 5).Value = '  -1.66%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''9.56'
$ws.Cells.Item(22, 5).Value = '  +3.24%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''0.728'
$ws.Cells.Item(23, 5).Value = '  -2.05%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''2.53'
$ws.Cells.Item(24, 5).Value = '  +5.13%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''87.88'
$ws.Cells.Item(25, 5).Value = '  -0.08%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''13.21'
$ws.Cells.Item(26, 5).Value = '  -1.84%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''11.06'
$ws.Cells.Item(27, 5).Value = '  +2.87%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''0.0000139'
$ws.Cells.Item(28, 5).Value = '  +25.46%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''0.999'
$ws.Cells.Item(29, 5).Value = '  -0.07%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -1.93%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +3.54%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -4.53%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''31.46'
$ws.Cells.Item(33, 5).Value = '  -0.83%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''0.116'
$ws.Cells.Item(34, 5).Value = '  -0.33%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +0.09%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''6.22'
$ws.Cells.Item(36, 5).Value = '  +0.59%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +3.46%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'Kaspa'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(38, 4).Value = '''0.134'
$ws.Cells.Item(38, 5).Value = '  +4.84%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'TheGraph'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(39, 4).Value = '''0.335'
$ws.Cells.Item(39, 5).Value = '  -2.56%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''2.14'
$ws.Cells.Item(40, 5).Value = '  +1.07%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''51.09'
$ws.Cells.Item(41, 5).Value = '  +0.36%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''45.71'
$ws.Cells.Item(42, 5).Value = '  -1.89%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''8.75'
$ws.Cells.Item(43, 5).Value = '  -1.38%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''422.68'
$ws.Cells.Item(44, 5).Value = '  +7.58%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '3.045.15'
$ws.Cells.Item(45, 5).Value = '  -4.37%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''2.82'
$ws.Cells.Item(46, 5).Value = '  +1.04%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -0.89%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''27.63'
$ws.Cells.Item(48, 5).Value = '  -2.30%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Monero'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(49, 4).Value = '''138.72'
$ws.Cells.Item(49, 5).Value = '  +2.70%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'USDe'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(50, 4).Value = '''1.00'
$ws.Cells.Item(50, 5).Value = '  -0.03%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''2.48'
$ws.Cells.Item(51, 5).Value = '  +0.67%  '
